# КПЗ 1 Розклад занять.docx - add missing "Заняття"/"Здача" dates for
# group ПЗ-42 on the ЛР02 and ЛР03 rows of the schedule table.
#
# Table layout (1-based):
#   col1 = ПЗ-41 "Заняття" date   col2 = ПЗ-41 "Здача" date
#   col3 = ПЗ-42 "Заняття" date   col4 = ПЗ-42 "Здача" date
#   col5 = lesson code (Л../ПР../ЛР..)
#
# Row 8  -> ЛР02 (14.09 / 18.09 / <empty> / <empty>)
# Row 9  -> ЛР03 (20.09 / 24.09 / <empty> / <empty>)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rPrXml = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

function Set-DateCell($row, $col, $text) {
    $cell = $t.Cell($row, $col)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr>' + $rPrXml + '</w:pPr><w:r>' + $rPrXml + '<w:t>' + $text + '</w:t></w:r></w:p>'
    $cell.Range.InsertXML($xml)
}

# ЛР02 row: add the ПЗ-42 "Заняття" (21.09) and "Здача" (25.09) dates.
Set-DateCell 8 3 "21.09"
Set-DateCell 8 4 "25.09"

# ЛР03 row: add the ПЗ-42 "Заняття" (24.09) and "Здача" (27.09.) dates;
# the "Заняття" cell also gets a light-blue highlight.
Set-DateCell 9 3 "24.09"
$t.Cell(9, 3).Shading.BackgroundPatternColor = 15849926
Set-DateCell 9 4 "27.09."
